$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(33, 8).Value = 167.39394
$ws.Cells.Item(33, 9).Value = 111.44444
$ws.Cells.Item(33, 11).Value = 111.44444
$ws.Cells.Item(33, 13).Value = 117.55556
$ws.Cells.Item(45, 8).Value = 0
$ws.Cells.Item(45, 10).Value = 0
$ws.Cells.Item(45, 12).ClearContents()
$ws.Cells.Item(45, 14).Value = 0
$ws.Cells.Item(92, 8).Value = 993.5185
$ws.Cells.Item(92, 9).Value = 464.45
$ws.Cells.Item(92, 11).Value = 464.45
$ws.Cells.Item(92, 13).Value = 783.55
$ws.Cells.Item(94, 8).Value = 2480
$ws.Cells.Item(94, 9).Value = 2480
$ws.Cells.Item(94, 11).Value = 2480
$ws.Cells.Item(94, 13).Value = -2029
$ws.Cells.Item(96, 8).Value = 737.6667
$ws.Cells.Item(96, 9).Value = 513
$ws.Cells.Item(96, 10).Value = 850
$ws.Cells.Item(96, 11).Value = 1539
$ws.Cells.Item(96, 12).Value = 2550
$ws.Cells.Item(96, 13).Value = -166
$ws.Cells.Item(96, 14).Value = -5296
$ws.Cells.Item(97, 8).Value = 4133.3335
$ws.Cells.Item(97, 10).Value = 4133.3335
$ws.Cells.Item(97, 12).Value = 12400.0005
$ws.Cells.Item(97, 14).Value = -13392.0005
$ws.Cells.Item(100, 8).Value = 3000
$ws.Cells.Item(100, 9).Value = 3000
$ws.Cells.Item(100, 10).Value = 3000
$ws.Cells.Item(100, 11).Value = 3000
$ws.Cells.Item(100, 12).Value = 3000
$ws.Cells.Item(100, 13).Value = -2459
$ws.Cells.Item(100, 14).Value = -4082
$ws.Cells.Item(103, 8).Value = 577.7273
$ws.Cells.Item(103, 9).Value = 356.1111
$ws.Cells.Item(103, 10).Value = 1575
$ws.Cells.Item(103, 11).Value = 1068.3333
$ws.Cells.Item(103, 12).Value = 4725
$ws.Cells.Item(103, 13).Value = -482.3333
$ws.Cells.Item(103, 14).Value = -5897
$ws.Cells.Item(104, 8).Value = 759.5
$ws.Cells.Item(104, 9).Value = 922
$ws.Cells.Item(104, 10).Value = 272
$ws.Cells.Item(104, 11).Value = 2766
$ws.Cells.Item(104, 12).Value = 816
$ws.Cells.Item(104, 13).Value = -1019
$ws.Cells.Item(104, 14).Value = -4310
$ws.Cells.Item(113, 8).Value = 3150.7144
$ws.Cells.Item(113, 9).Value = 3115.3845
$ws.Cells.Item(113, 11).Value = 3115.3845
$ws.Cells.Item(113, 13).Value = 138.6154999999999
$ws.Cells.Item(137, 8).Value = 3847962.5
$ws.Cells.Item(137, 9).Value = 5883430.5
$ws.Cells.Item(137, 10).Value = 3188.889
$ws.Cells.Item(137, 11).Value = 17650291.5
$ws.Cells.Item(137, 12).Value = 9566.667000000001
$ws.Cells.Item(137, 13).Value = -17647741.5
$ws.Cells.Item(137, 14).Value = -14666.667

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(102, 8).Value = 71429816
$ws.Cells.Item(102, 9).Value = 142857140
$ws.Cells.Item(102, 10).Value = 2500
$ws.Cells.Item(102, 11).Value = 142857140
$ws.Cells.Item(102, 12).Value = 2500
$ws.Cells.Item(102, 13).Value = -142855518
$ws.Cells.Item(102, 14).Value = -5744
$ws.Cells.Item(122, 8).Value = 18520500
$ws.Cells.Item(122, 9).Value = 2472.5
$ws.Cells.Item(122, 11).Value = 7417.5
$ws.Cells.Item(122, 13).Value = -4967.5

$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(94, 8).Value = 698.9167
$ws.Cells.Item(94, 9).Value = 689.6667
$ws.Cells.Item(94, 10).Value = 726.6667
$ws.Cells.Item(94, 11).Value = 689.6667
$ws.Cells.Item(94, 12).Value = 726.6667
$ws.Cells.Item(94, 13).Value = -238.6667
$ws.Cells.Item(94, 14).Value = -1628.6667
$ws.Cells.Item(99, 8).Value = 1009.94446
$ws.Cells.Item(99, 9).Value = 1025.4546
$ws.Cells.Item(99, 11).Value = 1025.4546
$ws.Cells.Item(99, 13).Value = 472.5454
$ws.Cells.Item(105, 8).Value = 41668668
$ws.Cells.Item(105, 9).Value = 50001910
$ws.Cells.Item(105, 11).Value = 50001910
$ws.Cells.Item(105, 13).Value = -50000163

$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(31, 8).Value = 2032.5714
$ws.Cells.Item(31, 9).Value = 1267.6923
$ws.Cells.Item(31, 10).Value = 3275.5
$ws.Cells.Item(31, 11).Value = 1267.6923
$ws.Cells.Item(31, 12).Value = 3275.5
$ws.Cells.Item(31, 13).Value = -972.6922999999999
$ws.Cells.Item(31, 14).Value = -3865.5
$ws.Cells.Item(34, 8).Value = 2032.5714
$ws.Cells.Item(34, 9).Value = 1267.6923
$ws.Cells.Item(34, 10).Value = 3275.5
$ws.Cells.Item(34, 11).Value = 1267.6923
$ws.Cells.Item(34, 12).Value = 3275.5
$ws.Cells.Item(34, 13).Value = -1065.6923
$ws.Cells.Item(34, 14).Value = -3679.5
$ws.Cells.Item(105, 8).Value = 861.125
$ws.Cells.Item(105, 9).Value = 841.2857
$ws.Cells.Item(105, 11).Value = 841.2857
$ws.Cells.Item(105, 13).Value = 905.7143
$ws.Cells.Item(107, 8).Value = 439.4375
$ws.Cells.Item(107, 9).Value = 402.6
$ws.Cells.Item(107, 10).Value = 500.83334
$ws.Cells.Item(107, 11).Value = 402.6
$ws.Cells.Item(107, 12).Value = 500.83334
$ws.Cells.Item(107, 13).Value = 1517.4
$ws.Cells.Item(107, 14).Value = -4340.83334

$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(68, 8).Value = 1224.1852
$ws.Cells.Item(68, 9).Value = 1150
$ws.Cells.Item(68, 10).Value = 1304.0769
$ws.Cells.Item(68, 11).Value = 3450
$ws.Cells.Item(68, 12).Value = 3912.2307
$ws.Cells.Item(68, 13).Value = -2639
$ws.Cells.Item(68, 14).Value = -5534.2307
$ws.Cells.Item(71, 8).Value = 1224.1852
$ws.Cells.Item(71, 9).Value = 1150
$ws.Cells.Item(71, 10).Value = 1304.0769
$ws.Cells.Item(71, 11).Value = 10350
$ws.Cells.Item(71, 12).Value = 11736.6921
$ws.Cells.Item(71, 13).Value = -6294
$ws.Cells.Item(71, 14).Value = -19848.6921
$ws.Cells.Item(122, 8).Value = 881.3889
$ws.Cells.Item(122, 10).Value = 1553.5714
$ws.Cells.Item(122, 12).Value = 13982.1426
$ws.Cells.Item(122, 14).Value = -18882.1426
$ws.Cells.Item(123, 8).Value = 2486.6667
$ws.Cells.Item(123, 10).Value = 3298.3333
$ws.Cells.Item(123, 12).Value = 9894.999899999999
$ws.Cells.Item(123, 14).Value = -14794.9999
$ws.Cells.Item(124, 8).Value = 887.0714
$ws.Cells.Item(124, 9).Value = 607.25
$ws.Cells.Item(124, 10).Value = 999
$ws.Cells.Item(124, 11).Value = 1821.75
$ws.Cells.Item(124, 12).Value = 2997
$ws.Cells.Item(124, 13).Value = 3088.25
$ws.Cells.Item(124, 14).Value = -12817
$ws.Cells.Item(125, 8).Value = 3852.6316
$ws.Cells.Item(125, 9).Value = 1400
$ws.Cells.Item(125, 10).Value = 4728.5713
$ws.Cells.Item(125, 11).Value = 4200
$ws.Cells.Item(125, 12).Value = 14185.7139
$ws.Cells.Item(125, 13).Value = 720
$ws.Cells.Item(125, 14).Value = -24025.7139
$ws.Cells.Item(131, 8).Value = 942.01514
$ws.Cells.Item(131, 9).Value = 529.9
$ws.Cells.Item(131, 10).Value = 1015.6071
$ws.Cells.Item(131, 11).Value = 1589.7
$ws.Cells.Item(131, 12).Value = 3046.8213
$ws.Cells.Item(131, 13).Value = 3450.3
$ws.Cells.Item(131, 14).Value = -13126.8213

$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(80, 8).Value = 3668.2666
$ws.Cells.Item(80, 9).Value = 2350
$ws.Cells.Item(80, 10).Value = 3871.077
$ws.Cells.Item(80, 11).Value = 2350
$ws.Cells.Item(80, 12).Value = 3871.077
$ws.Cells.Item(80, 13).Value = -1352
$ws.Cells.Item(80, 14).Value = -5867.077
$ws.Cells.Item(83, 8).Value = 3668.2666
$ws.Cells.Item(83, 9).Value = 2350
$ws.Cells.Item(83, 10).Value = 3871.077
$ws.Cells.Item(83, 11).Value = 11750
$ws.Cells.Item(83, 12).Value = 19355.385
$ws.Cells.Item(83, 13).Value = -6758
$ws.Cells.Item(83, 14).Value = -29339.385
$ws.Cells.Item(97, 8).Value = 1472.4117
$ws.Cells.Item(97, 9).Value = 1725.8334
$ws.Cells.Item(97, 10).Value = 864.2
$ws.Cells.Item(97, 11).Value = 1725.8334
$ws.Cells.Item(97, 12).Value = 864.2
$ws.Cells.Item(97, 13).Value = -1229.8334
$ws.Cells.Item(97, 14).Value = -1856.2
$ws.Cells.Item(102, 8).Value = 1799.8
$ws.Cells.Item(102, 9).Value = 1799.8
$ws.Cells.Item(102, 11).Value = 1799.8
$ws.Cells.Item(102, 13).Value = -177.8
$ws.Cells.Item(126, 8).Value = 1701.5555
$ws.Cells.Item(126, 9).Value = 1314.2858
$ws.Cells.Item(126, 10).Value = 3057
$ws.Cells.Item(126, 11).Value = 3942.8574
$ws.Cells.Item(126, 12).Value = 9171
$ws.Cells.Item(126, 13).Value = -1472.8574
$ws.Cells.Item(126, 14).Value = -14111
$ws.Cells.Item(132, 8).Value = 113688.555
$ws.Cells.Item(132, 9).Value = 101852.6
$ws.Cells.Item(132, 10).Value = 128483.5
$ws.Cells.Item(132, 11).Value = 305557.8
$ws.Cells.Item(132, 12).Value = 385450.5
$ws.Cells.Item(132, 13).Value = -303027.8
$ws.Cells.Item(132, 14).Value = -390510.5

$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(7, 8).Value = 1822.7778
$ws.Cells.Item(7, 9).Value = 1842.8572
$ws.Cells.Item(7, 10).Value = 1752.5
$ws.Cells.Item(7, 11).Value = 1842.8572
$ws.Cells.Item(7, 12).Value = 1752.5
$ws.Cells.Item(7, 13).Value = -1730.8572
$ws.Cells.Item(7, 14).Value = -1976.5
$ws.Cells.Item(61, 8).Value = 2481.3333
$ws.Cells.Item(61, 9).Value = 2282.2354
$ws.Cells.Item(61, 11).Value = 2282.2354
$ws.Cells.Item(61, 13).Value = -2080.2354
$ws.Cells.Item(113, 8).Value = 2481.3333
$ws.Cells.Item(113, 9).Value = 2282.2354
$ws.Cells.Item(113, 11).Value = 2282.2354
$ws.Cells.Item(113, 13).Value = -112.2354
$ws.Cells.Item(122, 8).Value = 3008.9092
$ws.Cells.Item(122, 9).Value = 2999.7778
$ws.Cells.Item(122, 10).Value = 3050
$ws.Cells.Item(122, 11).Value = 8999.3334
$ws.Cells.Item(122, 12).Value = 9150
$ws.Cells.Item(122, 13).Value = -6549.3334
$ws.Cells.Item(122, 14).Value = -14050
$ws.Cells.Item(126, 8).Value = 1822.7778
$ws.Cells.Item(126, 9).Value = 1842.8572
$ws.Cells.Item(126, 10).Value = 1752.5
$ws.Cells.Item(126, 11).Value = 5528.571599999999
$ws.Cells.Item(126, 12).Value = 5257.5
$ws.Cells.Item(126, 13).Value = -3058.571599999999
$ws.Cells.Item(126, 14).Value = -10197.5

$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(96, 8).Value = 6141.75
$ws.Cells.Item(96, 9).Value = 2040.6
$ws.Cells.Item(96, 10).Value = 9071.143
$ws.Cells.Item(96, 11).Value = 2040.6
$ws.Cells.Item(96, 12).Value = 9071.143
$ws.Cells.Item(96, 13).Value = -667.5999999999999
$ws.Cells.Item(96, 14).Value = -11817.143
$ws.Cells.Item(100, 8).Value = 111562.11
$ws.Cells.Item(100, 9).Value = 62944.875
$ws.Cells.Item(100, 10).Value = 500500
$ws.Cells.Item(100, 11).Value = 125889.75
$ws.Cells.Item(100, 12).Value = 1001000
$ws.Cells.Item(100, 13).Value = -125348.75
$ws.Cells.Item(100, 14).Value = -1002082
$ws.Cells.Item(122, 8).Value = 3349.3
$ws.Cells.Item(122, 9).Value = 1477
$ws.Cells.Item(122, 10).Value = 4597.5
$ws.Cells.Item(122, 11).Value = 4431
$ws.Cells.Item(122, 12).Value = 13792.5
$ws.Cells.Item(122, 13).Value = -1981
$ws.Cells.Item(122, 14).Value = -18692.5
$ws.Cells.Item(126, 8).Value = 1261.625
$ws.Cells.Item(126, 9).Value = 1227.5714
$ws.Cells.Item(126, 10).Value = 1500
$ws.Cells.Item(126, 11).Value = 3682.7142
$ws.Cells.Item(126, 12).Value = 4500
$ws.Cells.Item(126, 13).Value = -1212.7142
$ws.Cells.Item(126, 14).Value = -9440
